$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet name ---
$ws.Name = "Planilha1"

# --- Column widths (calibrated so the stored OOXML "width" matches the target values) ---
$ws.Range("A:A").ColumnWidth = 21.1668
$ws.Range("B:B").ColumnWidth = 21.8335
$ws.Range("C:C").ColumnWidth = 22.1668
$ws.Range("D:D").ColumnWidth = 19.8335
$ws.Range("E:E").ColumnWidth = 32.1668

# --- Row heights ---
$ws.Range("A1:E1").RowHeight = 15
$ws.Range("A2:E2").RowHeight = 15
$ws.Range("A3:E3").RowHeight = 15
$ws.Range("A4:E4").RowHeight = 15
$ws.Range("A5:E5").RowHeight = 15
$ws.Range("A6:E6").RowHeight = 15
$ws.Range("A7:E7").RowHeight = 15
$ws.Range("A8:E8").RowHeight = 15
$ws.Range("A9:E9").RowHeight = 15
$ws.Range("A10:E10").RowHeight = 15

# --- Header row values (row 1) ---
$ws.Range("A1").Value = "Fundos"
$ws.Range("B1").Value = "Topos"
$ws.Range("C1").Value = "Timestamp"
$ws.Range("D1").Value = "Nome da acao"
$ws.Range("E1").Value = "Valor atual"

# --- Row 2 values ---
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 500
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2020-05-14T23:20:56.592Z"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "WINM20"
$ws.Range("E2").Value = 130.84

# --- Row 3 values ---
$ws.Range("A3").Value = 100
$ws.Range("B3").Value = 600
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2020-05-15T23:20:56.592Z"

# --- Row 4 values ---
$ws.Range("A4").Value = 300
$ws.Range("B4").Value = 700
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2020-05-16T23:20:56.592Z"

# --- Row 5 values ---
$ws.Range("A5").Value = 400
$ws.Range("B5").Value = 800
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2020-05-17T23:20:56.592Z"
